$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.091.72"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.749.48"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.01"
$ws.Range("E5").Value = "  +4.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5281"
$ws.Range("E7").Value = "  +2.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2791"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "1.744.08"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07181"
$ws.Range("E11").Value = "  +2.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.40"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6457"
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.627"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "78.50"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9996"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9999"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "25.995.35"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006727"
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("D21").Value = "1.968.49"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.312"
$ws.Range("E22").Value = "  +5.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.742"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.239"
$ws.Range("E24").Value = "  +2.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.91"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.510"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.29"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.803"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "104.48"
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08282"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.798"
$ws.Range("E31").Value = "  +5.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.654"
$ws.Range("E32").Value = "  +7.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04568"
$ws.Range("E33").Value = "  +3.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.644"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.005"
$ws.Range("E35").Value = "  +3.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6348"
$ws.Range("E36").Value = "  +6.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01597"
$ws.Range("E38").Value = "  +3.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.947"
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9994"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.57"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3927"
$ws.Range("E42").Value = "  +2.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7451"
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.025"
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1146"
$ws.Range("E45").Value = "  +3.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.349"
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05351"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.94"
$ws.Range("E48").Value = "  +4.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.23"
$ws.Range("E49").Value = "  +4.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.615"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3450"
$ws.Range("E51").Value = "  +2.05%  "
